$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "56.702.05"
$ws.Range("E2").Value = "  +0.36%  "
# Row 3
$ws.Range("D3").Value = "2.387.76"
$ws.Range("E3").Value = "  +0.63%  "
# Row 4
$ws.Range("E4").Value = "  -0.31%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "503.59"
$ws.Range("E5").Value = "  -0.24%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "132.17"
$ws.Range("E6").Value = "  +1.92%  "
# Row 7
$ws.Range("E7").Value = "  -0.15%  "
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.550"
$ws.Range("E8").Value = "  -0.55%  "
# Row 9
$ws.Range("D9").Value = "2.393.67"
$ws.Range("E9").Value = "  -0.05%  "
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0971"
$ws.Range("E10").Value = "  +0.83%  "
# Row 11
$ws.Range("E11").Value = "  -0.51%  "
# Row 12
$ws.Range("E12").Value = "  +0.41%  "
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.65"
$ws.Range("E13").Value = "  +0.03%  "
# Row 14
$ws.Range("D14").Value = "2.814.09"
$ws.Range("E14").Value = "  +0.50%  "
# Row 15
$ws.Range("D15").Value = "56.629.09"
$ws.Range("E15").Value = "  +0.44%  "
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.63"
$ws.Range("E16").Value = "  +0.28%  "
# Row 17
$ws.Range("E17").Value = "  +0.81%  "
# Row 18
$ws.Range("D18").Value = "2.394.86"
$ws.Range("E18").Value = "  +0.25%  "
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.21"
$ws.Range("E19").Value = "  -0.13%  "
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.05"
$ws.Range("E20").Value = "  +0.10%  "
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "309.36"
$ws.Range("E21").Value = "  -0.86%  "
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.27"
$ws.Range("E22").Value = "  +0.76%  "
# Row 23
$ws.Range("E23").Value = "  +0.13%  "
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.57"
$ws.Range("E24").Value = "  -4.81%  "
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "67.07"
$ws.Range("E25").Value = "  +2.31%  "
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  -0.04%  "
# Row 27
$ws.Range("E27").Value = "  +0.83%  "
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.150"
$ws.Range("E28").Value = "  -0.50%  "
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.41"
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "176.25"
$ws.Range("E30").Value = "  +0.82%  "
# Row 31
$ws.Range("D31").Value = "0.0₃0723"
$ws.Range("E31").Value = "  +1.48%  "
# Row 32
$ws.Range("E32").Value = "  -1.05%  "
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.12"
$ws.Range("E33").Value = "  +1.15%  "
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.84"
$ws.Range("E34").Value = "  -4.78%  "
# Row 35
$ws.Range("E35").Value = "  +0.03%  "
# Row 36
$ws.Range("E36").Value = "  +0.30%  "
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "17.82"
$ws.Range("E37").Value = "  +0.21%  "
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.19"
$ws.Range("E38").Value = "  -2.09%  "
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.81"
$ws.Range("E39").Value = "  +1.09%  "
# Row 40
$ws.Range("B40").Value = "OKB"
$ws.Range("C40").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "36.85"
$ws.Range("E40").Value = "  +2.83%  "
# Row 41
$ws.Range("B41").Value = "SuiNetwork"
$ws.Range("C41").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.824"
$ws.Range("E41").Value = "  +5.00%  "
# Row 42
$ws.Range("E42").Value = "  +0.42%  "
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "131.39"
$ws.Range("E43").Value = "  +0.22%  "
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.38"
$ws.Range("E44").Value = "  +0.44%  "
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.86"
$ws.Range("E45").Value = "  +1.48%  "
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.567"
$ws.Range("E46").Value = "  -0.47%  "
# Row 47
$ws.Range("B47").Value = "Stellar"
$ws.Range("C47").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0909"
$ws.Range("E47").Value = "  +0.97%  "
# Row 48
$ws.Range("B48").Value = "Bittensor"
$ws.Range("C48").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "249.80"
$ws.Range("E48").Value = "  -2.11%  "
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0483"
$ws.Range("E49").Value = "  -0.95%  "
# Row 50
$ws.Range("E50").Value = "  +0.98%  "
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "17.18"
$ws.Range("E51").Value = "  +7.95%  "
